$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the workbook window (best-effort; mirrors the author's window resize)
$excel.ActiveWindow.Width = 14420

# Row 9: coverage simulation entry for 2023-08-03 (serial date 45141)
# Copy the date formatting from B8 (style used by the other date cells)
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = 45141

# Daily total (D) and the day-over-day delta formula (C), matching row 8's pattern
$ws.Range("D9").Value = 680
$ws.Range("C9").Formula = "=D9-D8"

# Move the active selection to D11 (where the author continued writing/entering data)
$ws.Range("D11").Select() | Out-Null
